$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "Login with valid username and password"
$ws.Range("B54").Value = "PASSED"
$ws.Range("C54").Value = "chrome"
$ws.Range("D54").Value = "11_06_23_152922"

$ws.Range("A55").Value = "Create Country"
$ws.Range("B55").Value = "PASSED"
$ws.Range("C55").Value = "chrome"
$ws.Range("D55").Value = "11_06_23_152931"
